$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cells that just get their number replaced.
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "6504"
$t.Cell(6,1).Range.Text  = "0.28854"
$t.Cell(7,1).Range.Text  = "0.08496"
$t.Cell(8,1).Range.Text  = "0.00572"
$t.Cell(9,1).Range.Text  = "0.24311"
$t.Cell(10,1).Range.Text = "0.25040"
$t.Cell(11,1).Range.Text = "0.26548"
$t.Cell(12,1).Range.Text = "24.02551"

# Rows 44-46 collapse their tab-separated multi-run content down to a
# single value each (moved from the old rows 1-3).
$t.Cell(44,1).Range.Text = "98.58"
$t.Cell(45,1).Range.Text = "24.03"
$t.Cell(46,1).Range.Text = "1687"
